$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new skill entry (row 12) below the existing data.
# Inserting the row (rather than just writing into empty cells) makes
# Excel carry over the formatting of the row above it, same as a user
# pressing "Insert Sheet Row" / dragging the table down one row.
$newRow = $ws.Rows("12:12")
$newRow.Insert(-4121, 0) | Out-Null

$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "CureI"
$ws.Range("C12").Value = "StatusSkill"
$ws.Range("D12").Value = 8
$ws.Range("E12").Value = 3

# Select the freshly added row, matching the last interaction recorded
# in the sheet (clicking the row-12 header).
$ws.Rows("12:12").Select() | Out-Null
